$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: extend the year header series with 2019 / 2020, matching the
#     formatting already used by the existing D3:I3 header cells (style 12).
$ws.Range("I3").Copy()
$ws.Range("J3:K3").PasteSpecial(-4122)
$ws.Range("J3").Value = 2019
$ws.Range("K3").Value = 2020

# --- Row 4: extend the data series with two more values (6.18, 6.18),
#     matching the formatting already used by D4:I4 (font/border/alignment).
$ws.Range("I4").Copy()
$ws.Range("J4:K4").PasteSpecial(-4122)
$ws.Range("J4").Value = 6.18
$ws.Range("K4").Value = 6.18

# Touch the interior/fill of the new cells so the style is recorded with an
# explicit (no-op) fill application, same as Excel leaves behind after a
# format round-trip -- matches the new cellXfs entry added by the edit.
$ws.Range("J4:K4").Interior.ColorIndex = 1
$ws.Range("J4:K4").Interior.Pattern = -4142

# --- Update the active selection left behind by the edit.
$ws.Range("G11").Select()
